# Append 45 new data rows (rows 102-146) to the
# "master-reg_center_device_h" master-data sheet, mirroring the pattern
# already used by the existing rows (regcntr_id cycling 10002-10010,
# device_id incrementing sequentially, and the remaining columns holding
# the same constant lookup values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (regcntr_id, device_id) for rows 102..146, in order.
$data = @(
    @(10002,3000121),
    @(10003,3000122),
    @(10004,3000123),
    @(10005,3000124),
    @(10006,3000125),
    @(10007,3000126),
    @(10008,3000127),
    @(10009,3000128),
    @(10010,3000129),
    @(10002,3000130),
    @(10003,3000131),
    @(10004,3000132),
    @(10005,3000133),
    @(10006,3000134),
    @(10007,3000135),
    @(10008,3000136),
    @(10009,3000137),
    @(10010,3000138),
    @(10002,3000139),
    @(10003,3000140),
    @(10004,3000141),
    @(10005,3000142),
    @(10006,3000143),
    @(10007,3000144),
    @(10008,3000145),
    @(10009,3000146),
    @(10010,3000147),
    @(10002,3000148),
    @(10003,3000149),
    @(10004,3000150),
    @(10005,3000151),
    @(10006,3000152),
    @(10007,3000153),
    @(10008,3000154),
    @(10009,3000155),
    @(10010,3000156),
    @(10002,3000157),
    @(10003,3000158),
    @(10004,3000159),
    @(10005,3000160),
    @(10006,3000161),
    @(10007,3000162),
    @(10008,3000163),
    @(10009,3000164),
    @(10010,3000165)
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $regcntrId = $data[$i][0]
    $deviceId = $data[$i][1]

    $ws.Cells.Item($row, 1).Value = $regcntrId      # A: regcntr_id
    $ws.Cells.Item($row, 2).Value = $deviceId       # B: device_id
    $ws.Cells.Item($row, 3).Value = "eng"           # C: lang_code
    $ws.Cells.Item($row, 4).Value = $true            # D: is_active
    $ws.Cells.Item($row, 5).Value = "superadmin"    # E: cr_by
    $ws.Cells.Item($row, 6).Value = "now()"         # F: cr_dtimes
    $ws.Cells.Item($row, 7).Value = "now()"         # G: eff_dtimes
}

# Reflect what was selected/visible in the workbook after the paste -
# the user had highlighted the newly added regcntr_id/device_id columns.
$lastRow = $startRow + $data.Count - 1
$ws.Range("A$($startRow):B$($lastRow)").Select()

# Make the printed page orientation explicit (matches the saved pageSetup).
$ws.PageSetup.Orientation = 1
